# Superficie_pop_2020.xlsx: the "Superficie (km² )" / "Population (habitants)"
# columns (B:C) on rows 2-19 currently hold the figures as TEXT (the column is
# formatted "@"), each stored as its own shared string (e.g. "69 711"). The
# authored edit retypes that block as genuine numeric data (e.g. 69711) while
# keeping every other formatting attribute (cell style / "@" text number
# format, alignment, etc.) exactly as it was, and updates the saved
# selection to the now-populated B2:C19 block.
#
# A plain `Range.Value = <number>` on a "@"-formatted cell gets re-coerced
# back to text by Excel (matches real Excel behaviour), so instead each cell
# is reset to the sheet's default/general style, given its numeric value,
# and then has the original "@" text formatting pasted back on top of it
# (format-only paste) so the stored style index is reused unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

$data = @(
    @{ Row = 2;  Superficie = 69711;  Population = 8064146 },
    @{ Row = 3;  Superficie = 47784;  Population = 2794517 },
    @{ Row = 4;  Superficie = 27208;  Population = 3358524 },
    @{ Row = 5;  Superficie = 39151;  Population = 2565726 },
    @{ Row = 6;  Superficie = 8680;   Population = 345867 },
    @{ Row = 7;  Superficie = 57433;  Population = 5536002 },
    @{ Row = 8;  Superficie = 1703;   Population = 379707 },
    @{ Row = 9;  Superficie = 83534;  Population = 288086 },
    @{ Row = 10; Superficie = 31813;  Population = 5987795 },
    @{ Row = 11; Superficie = 12011;  Population = 12291557 },
    @{ Row = 12; Superficie = 1128;   Population = 359821 },
    @{ Row = 13; Superficie = 2504;   Population = 278926 },
    @{ Row = 14; Superficie = 29906;  Population = 3313432 },
    @{ Row = 15; Superficie = 83809;  Population = 6018424 },
    @{ Row = 16; Superficie = 72724;  Population = 5951850 },
    @{ Row = 17; Superficie = 32082;  Population = 3818421 },
    @{ Row = 18; Superficie = 2505;   Population = 856858 },
    @{ Row = 19; Superficie = 31400;  Population = 5077583 }
)

# A never-edited cell whose style is the sheet's implicit default (General
# number format, no special alignment/font) - used to neutralise a cell's
# formatting before writing a fresh numeric literal into it.
$defaultStyleCell = $ws.Range("A21")

# Snapshot the data columns' original "@" text formatting (style index 2)
# onto an unused scratch cell before any edits touch B2:C19, so the exact
# formatting can be reapplied afterwards without Excel minting a new style
# table entry.
$scratch = $ws.Range("Z500")
$ws.Range("B2").Copy()
$scratch.PasteSpecial($xlPasteFormats)

foreach ($item in $data) {
    $r = $item.Row

    $bCell = $ws.Cells.Item($r, 2)
    $bCell.Style = $defaultStyleCell.Style
    $bCell.Value = $item.Superficie

    $cCell = $ws.Cells.Item($r, 3)
    $cCell.Style = $defaultStyleCell.Style
    $cCell.Value = $item.Population
}

# Re-apply the original "@" text formatting (without touching the values
# just written) to every cell in the block.
$scratch.Copy()
foreach ($item in $data) {
    $r = $item.Row
    $ws.Cells.Item($r, 2).PasteSpecial($xlPasteFormats)
    $ws.Cells.Item($r, 3).PasteSpecial($xlPasteFormats)
}

$scratch.Clear()

# Match the authored selection state.
$ws.Range("B2:C19").Select()
